# Generate Report for Handback
# Updates the "Latest Handback DateTime" (column K) for the row of the
# d19a2859-f3bf-46b5-97d4-0466a0e4744c.md file (row 2) on both the
# zh-cn and de-de localization-status sheets, reflecting a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("K2").Value = "2016-09-06 05:23:27"
$dede.Range("K2").Value = "2016-09-06 05:23:45"
